$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '47.639.53'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.494.69'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.91'
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.32'
$ws.Range("E6").Value = '  +3.82%  '
$ws.Range("E7").Value = '  -0.51%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.57'
$ws.Range("E10").Value = '  +4.13%  '
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.63'
$ws.Range("E13").Value = '  +2.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.21'
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.883.90'
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.494.66'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '47.505.68'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.42'
$ws.Range("E19").Value = '  +5.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.65'
$ws.Range("E20").Value = '  +1.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0943'
$ws.Range("E21").Value = '  +0.77%  '
$ws.Range("E22").Value = '  +15.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.71'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '247.08'
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.80'
$ws.Range("E27").Value = '  -1.31%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.00'
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("E30").Value = '  +3.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.85'
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("E32").Value = '  +0.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.41'
$ws.Range("E33").Value = '  +2.85%  '
$ws.Range("E34").Value = '  -0.41%  '
$ws.Range("E35").Value = '  +0.94%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.73'
$ws.Range("E37").Value = '  +2.85%  '
$ws.Range("E38").Value = '  +1.16%  '
$ws.Range("E39").Value = '  -1.35%  '
$ws.Range("E40").Value = '  +0.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.47'
$ws.Range("E41").Value = '  +6.10%  '
$ws.Range("E42").Value = '  -1.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '118.95'
$ws.Range("E43").Value = '  -1.84%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.998.56'
$ws.Range("E45").Value = '  +1.86%  '
$ws.Range("E46").Value = '  +2.38%  '
$ws.Range("E47").Value = '  -2.59%  '
$ws.Range("E48").Value = '  +0.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.08'
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.23'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '56.75'
$ws.Range("E51").Value = '  +3.46%  '
